# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" on the Overview sheet and the
# corresponding "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the per-language detail sheets for the
# 82dff81f-5794-4f20-8e52-1001d21d352b.md file, reflecting a freshly
# generated handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 2 is 82dff81f-5794-4f20-8e52-1001d21d352b.md
# Column G = "Latest HO Xliff Generate Date"
$overview.Range("G2").Value = "2016-08-27 06:45:20"

# zh-cn detail sheet: row 2 is 82dff81f-5794-4f20-8e52-1001d21d352b.md
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$zhcn.Range("H2").Value = "2016-08-27 06:45:16"
$zhcn.Range("K2").Value = "2016-08-27 06:45:32"

# de-de detail sheet: row 2 is 82dff81f-5794-4f20-8e52-1001d21d352b.md
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$dede.Range("H2").Value = "2016-08-27 06:45:20"
$dede.Range("K2").Value = "2016-08-27 06:45:39"
